$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '31.397.08'
$ws.Range('E2').Value = '  +3.30%  '
$ws.Range('D3').Value = '2.003.45'
$ws.Range('E3').Value = '  +7.03%  '
$c = $ws.Range('D4')
$c.Value = "'0.9968"
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.30%  '
$c = $ws.Range('D5')
$c.Value = "'0.8141"
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +72.76%  '
$c = $ws.Range('D6')
$c.Value = "'255.58"
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +4.69%  '
$c = $ws.Range('D7')
$c.Value = "'0.9973"
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -0.26%  '
$c = $ws.Range('D8')
$c.Value = "'0.3585"
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +24.15%  '
$c = $ws.Range('D9')
$c.Value = "'25.86"
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +17.53%  '
$c = $ws.Range('D10')
$c.Value = "'0.07016"
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +7.85%  '
$c = $ws.Range('D11')
$c.Value = "'0.8479"
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +16.92%  '
$c = $ws.Range('D12')
$c.Value = "'0.08143"
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +4.53%  '
$c = $ws.Range('D13')
$c.Value = "'101.64"
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +5.48%  '
$ws.Range('D14').Value = '1.993.34'
$ws.Range('E14').Value = '  +6.55%  '
$c = $ws.Range('D15')
$c.Value = "'5.539"
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +7.54%  '
$c = $ws.Range('D16')
$c.Value = "'272.59"
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -3.70%  '
$ws.Range('D17').Value = '31.369.86'
$ws.Range('E17').Value = '  +3.25%  '
$c = $ws.Range('D18')
$c.Value = "'14.04"
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +7.71%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D19')
$c.Value = "'5.845"
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +10.94%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range('D20')
$c.Value = "'0.000007959"
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +5.75%  '
$ws.Range('D21').Value = '2.251.11'
$ws.Range('E21').Value = '  +6.62%  '
$c = $ws.Range('D22')
$c.Value = "'0.9984"
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.15%  '
$c = $ws.Range('D23')
$c.Value = "'0.9979"
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -0.21%  '
$c = $ws.Range('D24')
$c.Value = "'7.025"
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +12.46%  '
$c = $ws.Range('D25')
$c.Value = "'9.865"
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +8.70%  '
$c = $ws.Range('D26')
$c.Value = "'0.1529"
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +58.88%  '
$c = $ws.Range('D27')
$c.Value = "'164.10"
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +0.15%  '
$c = $ws.Range('D28')
$c.Value = "'20.18"
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +7.32%  '
$c = $ws.Range('D29')
$c.Value = "'2.256"
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +19.64%  '
$c = $ws.Range('D30')
$c.Value = "'1.576"
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +5.92%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D31')
$c.Value = "'4.615"
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +8.83%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D32')
$c.Value = "'1.358"
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +2.77%  '
$c = $ws.Range('D33')
$c.Value = "'4.365"
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +5.65%  '
$c = $ws.Range('D34')
$c.Value = "'0.05214"
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +7.71%  '
$c = $ws.Range('D35')
$c.Value = "'1.223"
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +8.95%  '
$c = $ws.Range('D36')
$c.Value = "'0.7643"
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +10.38%  '
$c = $ws.Range('D37')
$c.Value = "'2.756"
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +1.53%  '
$c = $ws.Range('D38')
$c.Value = "'0.02015"
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +6.22%  '
$ws.Range('E39').Value = '  +3.35%  '
$c = $ws.Range('D40')
$c.Value = "'6.656"
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +6.38%  '
$c = $ws.Range('D41')
$c.Value = "'0.4781"
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +12.82%  '
$c = $ws.Range('D42')
$c.Value = "'78.81"
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +3.55%  '
$c = $ws.Range('D43')
$c.Value = "'2.138"
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +10.34%  '
$c = $ws.Range('D44')
$c.Value = "'0.8619"
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +3.98%  '
$c = $ws.Range('D45')
$c.Value = "'104.69"
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +3.67%  '
$c = $ws.Range('D46')
$c.Value = "'0.9982"
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -0.07%  '
$c = $ws.Range('D47')
$c.Value = "'9.980"
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +3.17%  '
$c = $ws.Range('D48')
$c.Value = "'7.557"
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +8.15%  '
$c = $ws.Range('D49')
$c.Value = "'0.4414"
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +12.01%  '
$c = $ws.Range('D50')
$c.Value = "'36.99"
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +4.69%  '
$c = $ws.Range('D51')
$c.Value = "'0.1204"
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +13.63%  '
